# Update "想去人数" (interest count) figures for several cap events, and
# mark row 25 ("南昌·Cookie动漫嘉年华-赵路专场票") as sold out.
#
# These figures change identically on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# row -> new F value (count of people interested)
$fUpdates = @{
    2  = 130
    7  = 1192
    8  = 1493
    9  = 331
    10 = 369
    12 = 130
    16 = 265
    17 = 288
    19 = 1697
    20 = 63
    23 = 643
    25 = 329
    26 = 4069
    28 = 475
    29 = 251
    30 = 1055
    31 = 126
    33 = 386
    34 = 22
    35 = 177
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    # Row 25's ticket is now sold out - lowest price column switches from a
    # numeric price to the text "已售罄" (sold out).
    $ws.Range("G25").Value = "已售罄"
}
